$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a new worksheet "2022-Q4" right after "总计" and before "2022-Q3"
#    by duplicating the existing "2022-Q3" sheet (so it inherits the same
#    column layout / header row / formatting), then overwriting its data.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$q3Sheet = $wb.Worksheets.Item("2022-Q3")
$q3Sheet.Copy($q3Sheet, [System.Reflection.Missing]::Value)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q4"

# Extend the formatting of the last template row (row 12) down to the extra
# rows needed (13-16), since the new sheet needs 15 data rows (2022-Q3's
# template only has 11).
$newSheet.Range("A12:H12").Copy()
$newSheet.Range("A13:H16").PasteSpecial(-4122)

# Force columns B, D, E, F, G to be stored as text so numeric-looking
# strings (fund codes, percentages, ...) keep their original text form
# instead of being auto-converted to numbers.
$newSheet.Range("B2:B16").NumberFormat = "@"
$newSheet.Range("D2:G16").NumberFormat = "@"

$q4Data = @(
    @(0, "008188", "前海开源稳健增长三年持有期混合", "21.91", "91.82", "4.57", "1.0013", 10),
    @(1, "010826", "大成产业趋势混合A", "11.37", "93.99", "4.17", "0.4741", 7),
    @(2, "000690", "前海开源大海洋战略经济灵活配置混合", "4.63", "92.54", "6.36", "0.2945", 4),
    @(3, "001178", "前海开源再融资主题精选股票", "5.94", "92.55", "4.34", "0.2578", 10),
    @(4, "010296", "万家互联互通中国优势量化策略混合A", "4.37", "94.52", "5.69", "0.2487", 9),
    @(5, "010827", "大成产业趋势混合C", "3.42", "93.99", "4.17", "0.1426", 7),
    @(6, "011287", "前海开源聚慧三年持有期混合", "2.88", "92.27", "4.48", "0.1290", 10),
    @(7, "006775", "前海开源优质成长混合", "2.55", "92.63", "4.60", "0.1173", 10),
    @(8, "003857", "前海开源周期优选灵活配置混合A", "2.13", "89.59", "5.42", "0.1154", 4),
    @(9, "000969", "前海开源大安全核心精选灵活配置混合", "1.09", "91.30", "5.35", "0.0583", 8),
    @(10, "006216", "前海开源价值成长灵活配置混合A", "1.14", "91.81", "4.44", "0.0506", 10),
    @(11, "003858", "前海开源周期优选灵活配置混合C", "0.72", "89.59", "5.42", "0.0390", 4),
    @(12, "010297", "万家互联互通中国优势量化策略混合C", "0.47", "94.52", "5.69", "0.0267", 9),
    @(13, "006217", "前海开源价值成长灵活配置混合C", "0.47", "91.81", "4.44", "0.0209", 10),
    @(14, "002020", "国都创新驱动灵活配置混合", "0.12", "83.47", "5.01", "0.0060", 2)
)

$r = 2
foreach ($row in $q4Data) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = $row[3]
    $newSheet.Cells.Item($r, 5).Value = $row[4]
    $newSheet.Cells.Item($r, 6).Value = $row[5]
    $newSheet.Cells.Item($r, 7).Value = $row[6]
    $newSheet.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: insert a new row for 2022-Q4 at the
#    top of the data (row 2), pushing the rest down by one row.
# ---------------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 15
$totalSheet.Range("D2").Value = 2.98

# Restore "总计" as the active sheet/tab (selection ends up on the last
# sheet we touched otherwise).
$totalSheet.Activate()
